$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Shape 1: "CuadroTexto 25" (id 26) ---------------------------------
$sh1 = $s.Shapes.Item(1)

# Reposition / resize the text box.
$sh1.Left   = 304.2274169921875
$sh1.Top    = 469.4731750488281
$sh1.Width  = 276.4996337890625
$sh1.Height = 138.1359100341797

$tr1 = $sh1.TextFrame.TextRange

# Paragraph 5: "Zonas de tiro( FGM, FGA, FGP)" -> "...FGPCT)"
$para5 = $tr1.Paragraphs(5, 1)
$para5.Runs(1, 1).Text = "Zonas de tiro( FGM, FGA, FGPCT)"

# Paragraph 6: merge the three runs ("Tiros " / "defencidos" / " (2pt, 3pt)")
# into a single corrected run "Tiros defendidos (2pt, 3pt)". Clear the
# trailing runs first (back to front) so the earlier run's growth doesn't
# shift the still-to-be-cleared runs out from under us, then grow run 1
# to hold the full corrected sentence (keeps its clean, non-"err" rPr).
$para6 = $tr1.Paragraphs(6, 1)
$para6.Runs(3, 1).Text = ""
$para6.Runs(2, 1).Text = ""
$para6.Runs(1, 1).Text = "Tiros defendidos (2pt, 3pt)"

# --- Shape 8: "Rectángulo: esquinas redondeadas 24" (id 25) ------------
$sh8 = $s.Shapes.Item(8)
$sh8.Left  = 292.20599365234375
$sh8.Width = 296.5939636230469

# --- Shape 15: "Conector: angular 54" (id 55) ---------------------------
$sh15 = $s.Shapes.Item(15)
$sh15.Left   = 429.9045104980469
$sh15.Top    = 427.14703369140625
$sh15.Height = 0.05417323112487793

# Mark the connector's non-visual properties as locked (adds <a:cxnSpLocks/>).
$sh15.LockAspectRatio = -1
